$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark from its current location (paragraph 1,
#    between the first and second tab runs).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Append a new paragraph at the end of the document that exercises
#    negative / hanging-indent tab-stop scenarios, re-creating the
#    _GoBack bookmark inside it (this is where Word leaves it after the
#    last edit made to the document).
$end = $d.Content
$end.Collapse(0)

$xmlFrag = '<?xml version="1.0" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p>' +
'<w:pPr>' +
'<w:tabs>' +
'<w:tab w:val="left" w:pos="-720"/>' +
'<w:tab w:val="left" w:pos="720"/>' +
'</w:tabs>' +
'<w:ind w:hanging="1080"/>' +
'</w:pPr>' +
'<w:r><w:tab/></w:r>' +
'<w:r><w:t>-0.5</w:t></w:r>' +
'<w:r><w:tab/></w:r>' +
'<w:r><w:tab/></w:r>' +
'<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
'<w:bookmarkEnd w:id="0"/>' +
'<w:r><w:t>0.5</w:t></w:r>' +
'</w:p>' +
'</w:body>' +
'</w:document>' +
'</pkg:xmlData>' +
'</pkg:part>' +
'</pkg:package>'

$end.InsertXML($xmlFrag)
